$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell B15: append a comma to the text "I'm Visuals"
$ws.Range("B15").Value = "I’m Visuals,"

# Move the active selection to D18 (reflected in sheet view selection)
$ws.Range("D18").Select()
